$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Finishing up cloth sim: log the 4 hours spent on "Kangassimulaatio" (row 21)
$ws.Range("G21").Value = 4

# Move the selection/cursor to F21 (and let the view scroll naturally,
# instead of staying pinned with A16 at the top).
$ws.Range("F21").Select()
